$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "Year of Treatment" column (column B), shifting all columns to its right one position to the left.
$ws.Range("B1").EntireColumn.Delete()

# Append ".global" suffix to each header label in row 1 (columns B through I after the shift).
$headerRange = $ws.Range("B1:I1")
for ($i = 1; $i -le $headerRange.Columns.Count; $i++) {
    $cell = $headerRange.Cells.Item(1, $i)
    $cell.Value = $cell.Value2 + ".global"
}
